# edit.ps1 — applies the "10.–16. 1. 2022" / "17.–23. 1. 2022" weekly-wave
# columns to both sheets of ZBP_06_home_office.xlsx, and bumps the
# "aktualizace" date in each sheet's trailing title row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": append columns BV (10.–16. 1. 2022) and BW (17.–23. 1. 2022)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Copy the style of the last existing header cell (BU1) onto the two new
# header cells so they keep the bold/centered/bordered look, then set text.
$ws1.Range("BU1").Copy()
$ws1.Range("BV1:BW1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("BV1").Value = "10.–16. 1. 2022"
$ws1.Range("BW1").Value = "17.–23. 1. 2022"

# Column BV = 74, column BW = 75
$s1_BV = @(0.76, 0.07000000000000001, 0.09, 0.08, 0.67, 0.08, 0.11, 0.14, 0.43, 0.12, 0.42, 0.03, 0.89, 0.04, 0.02, 0.05, 0.82, 0.05, 0.03, 0.1, 0.79, 0.16, 0.005, 0.045, 0.82, 0.04, 0.08, 0.06, 0.8100000000000001, 0.08, 0.04, 0.07000000000000001, 0.61, 0.12, 0.16, 0.11, 0.6, 0.12, 0.15, 0.13, 0.77, 0.08, 0.12, 0.03, 0.79, 0.06, 0.07000000000000001, 0.08, 0.63, 0.1, 0.09, 0.18, 0.74, 0.08, 0.07000000000000001, 0.11, 0.74, 0.12, 0.08, 0.06, 0.8, 0.03, 0.12, 0.05, 0.84, 0.02, 0.03, 0.11, 0.72, 0.11, 0.09, 0.08, 0.7, 0.09, 0.17, 0.04)
$s1_BW = @(0.73, 0.09, 0.09, 0.09, 0.67, 0.07000000000000001, 0.1, 0.16, 0.41, 0.14, 0.42, 0.03, 0.85, 0.07000000000000001, 0.02, 0.06, 0.78, 0.08, 0.06, 0.08, 0.78, 0.14, 0.02, 0.06, 0.8100000000000001, 0.06, 0.08, 0.05, 0.77, 0.08, 0.04, 0.11, 0.5600000000000001, 0.1, 0.15, 0.19, 0.57, 0.17, 0.15, 0.11, 0.76, 0.09, 0.12, 0.03, 0.75, 0.07000000000000001, 0.08, 0.1, 0.63, 0.14, 0.08, 0.15, 0.72, 0.09, 0.07000000000000001, 0.12, 0.72, 0.14, 0.07000000000000001, 0.07000000000000001, 0.77, 0.04, 0.13, 0.06, 0.82, 0.03, 0.03, 0.12, 0.7, 0.12, 0.1, 0.08, 0.66, 0.12, 0.17, 0.05)
for ($i = 0; $i -lt $s1_BV.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 74).Value = $s1_BV[$i]
    $ws1.Cells.Item($row, 75).Value = $s1_BW[$i]
}

# Update the "aktualizace" date in the trailing title row (row 78).
$ws1.Range("A78").Value = "Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 1. 2. 2022"

# ---------------------------------------------------------------------
# Sheet "pocetR": append columns BU (10.–16. 1. 2022) and BV (17.–23. 1. 2022)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

$ws2.Range("BT1").Copy()
$ws2.Range("BU1:BV1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("BU1").Value = "10.–16. 1. 2022"
$ws2.Range("BV1").Value = "17.–23. 1. 2022"

# Column BU = 73, column BV = 74
$s2_BU = @(1008, 282, 91, 262, 140, 87, 486, 248, 129, 145, 261, 566, 181, 459, 233, 316, 361, 400, 247)
$s2_BV = @(1008, 282, 91, 262, 140, 87, 486, 248, 129, 145, 261, 566, 181, 459, 233, 316, 361, 400, 247)
for ($i = 0; $i -lt $s2_BU.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 73).Value = $s2_BU[$i]
    $ws2.Cells.Item($row, 74).Value = $s2_BV[$i]
}

# Update the "aktualizace" date in the trailing title row (row 21), and
# extend the row's blank text-formatted trailer into the two new columns.
$ws2.Range("A21").Value = "Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 2. 2022"
$ws2.Cells.Item(21, 73).Font.Bold = $false
$ws2.Cells.Item(21, 74).Font.Bold = $false
